$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp on the existing row 4 (EMP-002)
$ws.Range("G4").Value = "2026-02-27 23:07:03"

# Copy formatting (style + row height) from row 4 into the two new rows
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)
$ws.Range("A4:H4").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 5: EMP-003
$ws.Range("A5").Value = "EMP-003"
$ws.Range("B5").Value = "Sujal Ashok Vaidya"
$ws.Range("C5").Value = "Manufacturing"
$ws.Range("D5").Value = "✗"
$ws.Range("E5").Value = "✗"
$ws.Range("F5").Value = "NOT READY"
$ws.Range("G5").Value = "2026-02-27 23:06:38"
$ws.Range("H5").Value = "Missing PPE: Helmet, Safety Vest"

# Row 6: EMP-004
$ws.Range("A6").Value = "EMP-004"
$ws.Range("B6").Value = "Vaibhav Hujaratti"
$ws.Range("C6").Value = "Electrical"
$ws.Range("D6").Value = "✗"
$ws.Range("E6").Value = "✗"
$ws.Range("F6").Value = "NOT READY"
$ws.Range("G6").Value = "2026-02-27 23:06:52"
$ws.Range("H6").Value = "Missing PPE: Helmet, Safety Vest"

# Ensure row heights match row 4 (22pt, custom height)
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight
$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(4).RowHeight
